$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 287
$ws1.Range("F4").Value = 3513
$ws1.Range("F5").Value = 2187
$ws1.Range("F8").Value = 67
$ws1.Range("F9").Value = 55
$ws1.Range("F10").Value = 1300
$ws1.Range("F12").Value = 1748
$ws1.Range("F13").Value = 132

# Sheet "全部类型" (All Types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 287
$ws4.Range("F4").Value = 3513
$ws4.Range("F5").Value = 2187
$ws4.Range("F9").Value = 67
$ws4.Range("F10").Value = 55
$ws4.Range("F13").Value = 1300
$ws4.Range("F15").Value = 1748
$ws4.Range("F16").Value = 132

$wb.Save()
